$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.528.69'
$ws.Range("E2").Value = '  +0.87%  '

$ws.Range("D3").Value = '3.169.19'
$ws.Range("E3").Value = '  -0.38%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '572.08'
$ws.Range("E5").Value = '  +0.32%  '

$ws.Range("D6").Value = '164.63'
$ws.Range("E6").Value = '  -2.60%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  -3.87%  '

$ws.Range("E9").Value = '  -2.68%  '

$ws.Range("D10").Value = '6.64'
$ws.Range("E10").Value = '  -0.99%  '

$ws.Range("D11").Value = '0.385'
$ws.Range("E11").Value = '  -0.48%  '

$ws.Range("D12").Value = '3.715.78'
$ws.Range("E12").Value = '  -0.66%  '

$ws.Range("E13").Value = '  -1.05%  '

$ws.Range("D14").Value = '64.520.48'
$ws.Range("E14").Value = '  +0.59%  '

$ws.Range("D15").Value = '25.40'
$ws.Range("E15").Value = '  -0.08%  '

$ws.Range("D16").Value = '3.163.22'
$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("E17").Value = '  -1.91%  '

$ws.Range("D18").Value = '409.54'
$ws.Range("E18").Value = '  -1.48%  '

$ws.Range("D19").Value = '12.80'
$ws.Range("E19").Value = '  -0.36%  '

$ws.Range("E20").Value = '  -1.71%  '

$ws.Range("D21").Value = '7.11'
$ws.Range("E21").Value = '  -0.81%  '

$ws.Range("E22").Value = '  +0.23%  '

$ws.Range("D23").Value = '68.86'
$ws.Range("E23").Value = '  -2.79%  '

$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").Value = '0.486'
$ws.Range("E24").Value = '  -1.65%  '

$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").Value = '0.197'
$ws.Range("E25").Value = '  -2.22%  '

$ws.Range("E26").Value = '  -5.99%  '

$ws.Range("D27").Value = '8.92'
$ws.Range("E27").Value = '  +1.44%  '

$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.36%  '

$ws.Range("D29").Value = '1.83'
$ws.Range("E29").Value = '  -1.46%  '

$ws.Range("D30").Value = '21.28'
$ws.Range("E30").Value = '  -3.03%  '

$ws.Range("D31").Value = '4.93'
$ws.Range("E31").Value = '  -1.58%  '

$ws.Range("D32").Value = '6.38'
$ws.Range("E32").Value = '  -0.52%  '

$ws.Range("E33").Value = '  -0.60%  '

$ws.Range("D34").Value = '156.20'
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("E35").Value = '  -2.04%  '

$ws.Range("E36").Value = '  -0.44%  '

$ws.Range("D37").Value = '2.690.08'
$ws.Range("E37").Value = '  -2.19%  '

$ws.Range("D38").Value = '24.04'
$ws.Range("E38").Value = '  -4.38%  '

$ws.Range("D39").Value = '4.11'
$ws.Range("E39").Value = '  -2.09%  '

$ws.Range("E40").Value = '  -2.82%  '

$ws.Range("D41").Value = '0.0622'
$ws.Range("E41").Value = '  -1.04%  '

$ws.Range("D42").Value = '5.45'
$ws.Range("E42").Value = '  -4.84%  '

$ws.Range("E43").Value = '  -1.83%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = '292.65'
$ws.Range("E44").Value = '  -1.81%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '21.54'
$ws.Range("E45").Value = '  -2.38%  '

$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").Value = '0.0988'
$ws.Range("E47").Value = '  -0.71%  '

$ws.Range("E48").Value = '  -7.21%  '

$ws.Range("D49").Value = '10.49'
$ws.Range("E49").Value = '  +0.70%  '

$ws.Range("D50").Value = '5.74'
$ws.Range("E50").Value = '  -1.37%  '
